$d = $word.ActiveDocument

# Helper: build the pPr XML shared by all the new "ListParagraph" numbered items.
function Get-ListPPr() {
    return '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
}

# Helper: insert a brand-new ListParagraph-styled, numbered paragraph right
# after $refPara, containing the runs described by $runsXml (already-built
# <w:r>...</w:r> markup). Returns the newly created Paragraph object so the
# caller can chain further insertions after it.
function Insert-ListParagraph($refPara, $runsXml) {
    $refPara.Range.InsertParagraphAfter()
    $newPara = $refPara.Next()

    $pPr = Get-ListPPr
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $pPr + $runsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $newPara.Range.InsertXML($xml)
    return $newPara
}

# Find the paragraph that ends with "Euclidean distance in Python?" -- the
# ten new questions are inserted right after it (and before the trailing
# empty list paragraph already present at the end of the document).
$ref = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Euclidean distance in Python\?") {
        $ref = $p
        break
    }
}

if ($ref -eq $null) {
    throw "Could not locate the 'Euclidean distance in Python?' paragraph"
}

$runs = '<w:r><w:t>What are dimensionality reduction and its benefits?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t>How will you calculate eigenvalues and eigenvectors of the following 3*3 matrix?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t>How should you maintain a deployed model?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t xml:space="preserve">What are recommender </w:t></w:r>' +
        '<w:r><w:t>systems?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t>How do you find RMSE and MSE in a linear regression model?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t>How can you select k for k-means?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t xml:space="preserve">What is the significance of </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">the </w:t></w:r>' +
        '<w:r><w:t>p-value?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t>How can outlier values be treated?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t>How can time series data be declared stationary?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs

$runs = '<w:r><w:t>How can you calculate accuracy using a confusion matrix?</w:t></w:r>'
$ref = Insert-ListParagraph $ref $runs
